# Printing all the graphs in one graph, calculating how long each child time travel
#
# The sheet tracks one "path" per row (one child per row). The commit resorts
# the children by their travel/wait time and adds a new column H with the
# computed travel duration (in minutes) for each child.
#
# Numeric-looking values in this sheet are stored as TEXT (the columns mix
# free-form strings like "7:00:00" / "-9.64,-6.26" with bare numbers like
# "14"), so plain-number cells are written with a leading apostrophe to force
# Excel to keep them as text instead of coercing them to the Number type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (child #0): swap in Nubia/Royce, update coords, contact and new duration
$ws.Range("B6").Value = "'15"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "Nubia  "
$ws.Range("D6").Value = "Royce  "
$ws.Range("E6").Value = "-9.16,-3.53"
$ws.Range("F6").Value = "Augustus(father): 0517389040"
$ws.Range("H6").Value = "'30.0"
$ws.Range("H6").Style = "Normal"

# Row 7 (child #1): swap in Elwanda/Cassy, update coords, contact and new duration
$ws.Range("B7").Value = "'2"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "Elwanda  "
$ws.Range("D7").Value = "Cassy  "
$ws.Range("E7").Value = "-7.38,-6.34"
$ws.Range("F7").Value = "Tamisha(mother): 0550693864"
$ws.Range("H7").Value = "'25.0"
$ws.Range("H7").Style = "Normal"

# Row 8 (child #2): swap in Fay/Emilee, update coords, contact, pickup time and duration
$ws.Range("B8").Value = "'13"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "Fay  "
$ws.Range("D8").Value = "Emilee  "
$ws.Range("E8").Value = "-3.16,-7.95"
$ws.Range("F8").Value = "Sheri(mother): 0516797453"
$ws.Range("G8").Value = "7:11:00"
$ws.Range("H8").Value = "'19.0"
$ws.Range("H8").Style = "Normal"

# Row 9 (child #3): swap in Lorinda/Tyron, update coords, contact, pickup time and duration
$ws.Range("B9").Value = "'14"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "Lorinda  "
$ws.Range("D9").Value = "Tyron  "
$ws.Range("E9").Value = "-0.99,-5.37"
$ws.Range("F9").Value = "Teresa(grandmother): 0558587699"
$ws.Range("G9").Value = "7:15:00"
$ws.Range("H9").Value = "'15.0"
$ws.Range("H9").Style = "Normal"

# Row 10 (child #4): swap in Trudie/Fleta, update coords, contact, pickup time and duration
$ws.Range("B10").Value = "'0"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "Trudie  "
$ws.Range("D10").Value = "Fleta  "
$ws.Range("E10").Value = "-4.17,-3.8"
$ws.Range("F10").Value = "Anneliese(father): 0548973345"
$ws.Range("G10").Value = "7:19:00"
$ws.Range("H10").Value = "'11.0"
$ws.Range("H10").Style = "Normal"

# Row 11 (child #5): swap in Wyatt/Willette, update coords, contact, pickup time and duration
$ws.Range("B11").Value = "'7"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "Wyatt  "
$ws.Range("D11").Value = "Willette  "
$ws.Range("E11").Value = "-2.87,-2.03"
$ws.Range("F11").Value = "Antionette(father): 0557331799"
$ws.Range("G11").Value = "7:23:00"
$ws.Range("H11").Value = "'7.0"
$ws.Range("H11").Style = "Normal"

# Row 12 (child #6): swap in Britta/Jamel, update coords, contact, pickup time and duration
$ws.Range("B12").Value = "'17"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "Britta  "
$ws.Range("D12").Value = "Jamel  "
$ws.Range("E12").Value = "-3.08,-1.04"
$ws.Range("F12").Value = "Albertine(father): 0574981040"
$ws.Range("G12").Value = "7:25:00"
$ws.Range("H12").Value = "'5.0"
$ws.Range("H12").Style = "Normal"

# Row 13 (school): arrival time moves later
$ws.Range("G13").Value = "7:30:00"

# Row 15: total trip time updates to match the new schedule
$ws.Range("B15").Value = "'30.0"
$ws.Range("B15").Style = "Normal"
